# Update the dSF column (F) with repulled/recalculated values.
# This mirrors a data "repull" where several dSF values differ from
# their previous (copy of dS0) values after recomputation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    5  = -3
    8  = -1
    12 = -2
    14 = -6
    20 = 5
    21 = -5
    22 = -4
    23 = -9
    24 = -6
    25 = -7
    27 = -7
    29 = -3
    34 = -4
    37 = -4
    38 = -3
    39 = -1
    40 = 9
    41 = 5
    42 = -2
    43 = 3
    45 = -1
    46 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
